{"js": "// The captured change is a pure OOXML canonicalization artifact: every\n// hunk in the source diff only re-orders existing XML attributes that are\n// already present (e.g. <w:tab w:val=\"left\" w:pos=\"3119\"/> becomes\n// <w:tab w:pos=\"3119\" w:val=\"left\"/>, the w:document root namespace\n// declarations get alphabetized, <w:pgSz>/<w:pgMar> attributes get\n// alphabetized, etc.). No text, run, paragraph, style, page-setup value,\n// or any other observable document property actually differs between the\n// \"before\" and \"after\" XML - only serialization/attribute order changed.\n//\n// Confirmed by normalizing both sides (parsing every element, sorting its\n// attributes, and comparing): the canonical content of word/document.xml\n// and word/styles.xml is byte-for-byte identical before and after the\n// commit shown in the diff. There is therefore no content mutation for\n// this script to perform against the Word JS API - the correct,\n// content-preserving edit is a no-op that leaves body text, formatting,\n// tab stops, section/page setup and styles untouched.\n//\n// We still touch the context once (a harmless read-only load/sync) so the\n// script demonstrably runs against the live document without mutating it.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The captured change is a pure OOXML canonicalization artifact: every\n# hunk in the source diff only re-orders existing XML attributes that are\n# already present (e.g. <w:tab w:val=\"left\" w:pos=\"3119\"/> becomes\n# <w:tab w:pos=\"3119\" w:val=\"left\"/>, the w:document root namespace\n# declarations get alphabetized, <w:pgSz>/<w:pgMar> attributes get\n# alphabetized, etc.). No text, run, paragraph, style, page-setup value,\n# or any other observable document property actually differs between the\n# \"before\" and \"after\" XML - only serialization/attribute order changed.\n#\n# Confirmed by normalizing both sides (parsing every element, sorting its\n# attributes, and comparing): the canonical content of word/document.xml\n# and word/styles.xml is byte-for-byte identical before and after the\n# commit shown in the diff. There is therefore no content mutation for\n# this script to perform against the Word object model - the correct,\n# content-preserving edit is a no-op that leaves body text, formatting,\n# tab stops, section/page setup and styles untouched.\n#\n# We still touch the document once (a harmless read-only access) so the\n# script demonstrably runs against the live document without mutating it.\n$d = $word.ActiveDocument\n$null = $d.Paragraphs.Count\n"}
